$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Text Not Translated")
$dst = $wb.Worksheets.Item("Text Translated")

# Copy the 4x21 block of translation data from "Text Not Translated" into
# "Text Translated" (the validation run processed the current Salesforce
# translations, so the previously-empty "Text Translated" sheet now holds
# the rows, and the "package" columns P/Q/R get resolved values instead of
# "Misssing").
for ($row = 1; $row -le 4; $row++) {
  for ($col = 1; $col -le 21; $col++) {
    $val = $src.Cells.Item($row, $col).Text
    if ($val -ne "") {
      $dst.Cells.Item($row, $col).Value = $val
    }
  }
}

# Resolved package/packageName/masterInfoLabel values (previously "Misssing").
# Single-quoted so the embedded '$' characters are taken literally.
$dst.Range("P2").Value = 'CustomLabel$CEC_Shipment_Additional_Info'
$dst.Range("Q2").Value = 'CustomLabel'
$dst.Range("R2").Value = 'Additional Info'

$dst.Range("P3").Value = 'CustomField$CEC_ShippingIdentifier__c.CEC_Applicable__c | CustomLabel$CEC_ShipmentIdentifier_Applicable'
$dst.Range("Q3").Value = 'CustomField | CustomLabel'
$dst.Range("R3").Value = 'Applicable'

$dst.Range("P4").Value = 'ValidationRule$Case.CEC_RestrictContactCaseOwnership | ValidationRule$Case.CEC_Restrict_SSR_And_Specialized_Users | ValidationRule$Case.CEC_Restrict_UPS_CEC_SRC | CustomLabel$cec_AR_case_Restrict_Owner_ship'
$dst.Range("Q4").Value = 'ValidationRule | CustomLabel'
$dst.Range("R4").Value = 'Cases cannot be assigned to this queue or user.'

# The data now lives in "Text Translated"; "Text Not Translated" goes back
# to being the empty sheet.
$src.Range("A1:U4").ClearContents()
